# Convert the roster time-in/time-out grid (rows 3-9) from 24-hour fraction-of-day
# values to 12-hour "h:mmAM/PM" text, matching the "change 24 hours time format to
# 12 hours format" commit. Also updates row numbers in column A, clears the blank
# cells in row 3 (G3/H3), resizes column E, and leaves the final selection on N3 -
# all exactly as captured in the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = 6
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "9:30AM"
$ws.Range("D3").Value = "12:10PM"
$ws.Range("E3").Value = "9:30AM"
$ws.Range("F3").Value = "12:00AM"
$ws.Range("G3").Value = ""
$ws.Range("H3").Value = ""
$ws.Range("I3").Value = "5:48PM"
$ws.Range("J3").Value = "5:40AM"
$ws.Range("K3").Value = "5:41AM"
$ws.Range("K3").NumberFormat = "h:mm"
$ws.Range("L3").Value = "5:42PM"
$ws.Range("L3").NumberFormat = "h:mm"
$ws.Range("M3").Value = "5:43AM"
$ws.Range("N3").Value = "5:44PM"

# Row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "12:30PM"
$ws.Range("D4").Value = "9:00AM"
$ws.Range("E4").Value = "9:31AM"
$ws.Range("F4").Value = "12:01AM"
$ws.Range("G4").Value = "9:40AM"
$ws.Range("H4").Value = "5:40PM"
$ws.Range("I4").Value = "5:49PM"
$ws.Range("J4").Value = "5:41AM"
$ws.Range("K4").Value = "5:42AM"
$ws.Range("L4").Value = "3:32PM"
$ws.Range("M4").Value = "5:44AM"
$ws.Range("N4").Value = "5:45AM"

# Row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "4:30AM"
$ws.Range("D5").Value = "12:30PM"
$ws.Range("E5").Value = "9:32AM"
$ws.Range("F5").Value = "12:02AM"
$ws.Range("G5").Value = "9:41AM"
$ws.Range("H5").Value = "5:41PM"
$ws.Range("I5").Value = "5:50PM"
$ws.Range("J5").Value = "5:42AM"
$ws.Range("K5").Value = "5:43AM"
$ws.Range("L5").Value = "3:33PM"
$ws.Range("M5").Value = "5:45AM"
$ws.Range("N5").Value = "5:46AM"

# Row 6
$ws.Range("A6").Value = 6
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = "5:30AM"
$ws.Range("D6").Value = "12:31PM"
$ws.Range("E6").Value = "9:33AM"
$ws.Range("F6").Value = "12:03AM"
$ws.Range("G6").Value = "9:42AM"
$ws.Range("H6").Value = "5:42PM"
$ws.Range("I6").Value = "5:51PM"
$ws.Range("J6").Value = "5:43AM"
$ws.Range("K6").Value = "5:44AM"
$ws.Range("L6").Value = "3:34PM"
$ws.Range("M6").Value = "5:46AM"
$ws.Range("N6").Value = "5:47AM"

# Row 7
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = "5:31AM"
$ws.Range("D7").Value = "12:32PM"
$ws.Range("E7").Value = "9:34AM"
$ws.Range("F7").Value = "12:04AM"
$ws.Range("G7").Value = "9:43AM"
$ws.Range("H7").Value = "5:43PM"
$ws.Range("I7").Value = "5:52PM"
$ws.Range("J7").Value = "5:44AM"
$ws.Range("K7").Value = "5:45AM"
$ws.Range("L7").Value = "3:35PM"
$ws.Range("M7").Value = "5:47AM"
$ws.Range("N7").Value = "5:48AM"

# Row 8
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = "5:32AM"
$ws.Range("D8").Value = "12:33PM"
$ws.Range("E8").Value = "9:35AM"
$ws.Range("F8").Value = "12:05AM"
$ws.Range("G8").Value = "9:44AM"
$ws.Range("H8").Value = "5:44PM"
$ws.Range("I8").Value = "5:53PM"
$ws.Range("J8").Value = "5:45AM"
$ws.Range("K8").Value = "5:46AM"
$ws.Range("L8").Value = "3:36PM"
$ws.Range("M8").Value = "5:78PM"
$ws.Range("N8").Value = "5:49AM"

# Row 9
$ws.Range("A9").Value = 6
$ws.Range("B9").Value = 3
$ws.Range("C9").Value = "5:33AM"
$ws.Range("D9").Value = "12:34PM"
$ws.Range("E9").Value = "9:36AM"
$ws.Range("F9").Value = "12:06AM"
$ws.Range("G9").Value = "9:45AM"
$ws.Range("H9").Value = "5:45PM"
$ws.Range("I9").Value = "5:54PM"
$ws.Range("J9").Value = "5:46AM"
$ws.Range("K9").Value = "5:47AM"
$ws.Range("L9").Value = "3:37PM"
$ws.Range("M9").Value = "5:79PM"
$ws.Range("N9").Value = "5:50AM"

# Column E was resized (narrow custom width) as part of the edit.
$ws.Columns("E:E").ColumnWidth = 9

# Final selection left on N3.
$ws.Range("N3").Select()
